$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Mark the four "practice" pairs (rows 2-5) as pair_kind = "generic" (col J) ---
foreach ($r in 2..5) {
    $ws.Range("J$r").Value = "generic"
}

# --- New "stim details" block added below the existing table ---
$ws.Range("A27").Value = "stim details"

# Header row for the new block
$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# Data rows for the new block: month, word_type pairs
$stimDetails = @(
    @(6, "video"),
    @(6, "video"),
    @(7, "video"),
    @(7, "video"),
    @(6, "audio"),
    @(6, "audio"),
    @(7, "audio"),
    @(7, "audio")
)

$row = 29
foreach ($entry in $stimDetails) {
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $row++
}
